# "Size scaling in same plot"
# The SA_series_scale_param (P) and V_series_scale_param (Q) columns are no
# longer needed as separate columns, so remove them entirely. Excel shifts
# the remaining columns (elastic_mod_mean, elastic_mod_std, elastic_mod_CV,
# log_elastic_mod_mean, cross_section) two places to the left to fill the
# gap, turning the old R:V range into the new P:T range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1:Q1").EntireColumn.Delete()
